# Updated cryptos list on Sat Dec  2 12:57:06 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}


# Row 2 - Bitcoin
Set-TextValue "D2" "38.773.44"
Set-TextValue "E2" "  +1.04%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.103.56"
Set-TextValue "E3" "  +0.83%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "227.62"
Set-TextValue "E5" "  -0.10%  "

# Row 6 - XRP
Set-TextValue "D6" "0.614"
Set-TextValue "E6" "  +0.33%  "

# Row 7 - Solana
Set-TextValue "D7" "62.43"
Set-TextValue "E7" "  +2.65%  "

# Row 8 - USDC
Set-TextValue "E8" "  -0.07%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +2.58%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +0.62%  "

# Row 11 - TRON
Set-TextValue "E11" "  -1.00%  "

# Row 12 - Chainlink
Set-TextValue "D12" "15.80"
Set-TextValue "E12" "  +6.77%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.416.28"
Set-TextValue "E13" "  +0.71%  "

# Row 14 - Avalanche
Set-TextValue "D14" "22.02"
Set-TextValue "E14" "  -0.96%  "

# Row 15 - Polygon
Set-TextValue "E15" "  +3.84%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.53"
Set-TextValue "E16" "  +1.76%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.086.08"
Set-TextValue "E17" "  +0.34%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "38.749.38"
Set-TextValue "E18" "  +1.16%  "

# Row 19 - Uniswap
Set-TextValue "D19" "6.13"
Set-TextValue "E19" "  +1.39%  "

# Row 20 - Litecoin
Set-TextValue "D20" "71.61"
Set-TextValue "E20" "  +0.98%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0843"
Set-TextValue "E21" "  +1.37%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "228.53"
Set-TextValue "E22" "  +1.50%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  -3.07%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +0.45%  "

# Row 26 - Monero (was Cosmos)
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "172.28"
Set-TextValue "E26" "  +1.69%  "

# Row 27 - Cosmos (was Monero)
Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "9.58"
Set-TextValue "E27" "  +1.91%  "

# Row 28 - Kaspa
Set-TextValue "E28" "  +2.81%  "

# Row 29 - ImmutableX
Set-TextValue "E29" "  +4.42%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "19.36"
Set-TextValue "E30" "  +1.93%  "

# Row 31 - WEMIXToken
Set-TextValue "D31" "2.55"
Set-TextValue "E31" "  +8.19%  "

# Row 32 - Stellar
Set-TextValue "E32" "  +0.76%  "

# Row 33 - Filecoin
Set-TextValue "D33" "4.56"
Set-TextValue "E33" "  +1.24%  "

# Row 34 - THORChain (was InternetComputer(DFINITY))
Set-TextValue "B34" "THORChain"
Set-TextValue "C34" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D34" "7.13"
Set-TextValue "E34" "  +11.52%  "

# Row 35 - InternetComputer(DFINITY) (was THORChain)
Set-TextValue "B35" "InternetComputer(DFINITY)"
Set-TextValue "C35" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D35" "4.75"
Set-TextValue "E35" "  -0.62%  "

# Row 36 - Hedera
Set-TextValue "E36" "  +1.96%  "

# Row 37 - LidoDAOToken
Set-TextValue "E37" "  +0.62%  "

# Row 38 - RenderToken
Set-TextValue "E38" "  -0.04%  "

# Row 39 - BinanceUSD
Set-TextValue "E39" "  -0.10%  "

# Row 40 - InjectiveProtocol
Set-TextValue "D40" "18.26"
Set-TextValue "E40" "  -0.53%  "

# Row 41 - Aave
Set-TextValue "D41" "102.79"
Set-TextValue "E41" "  +2.99%  "

# Row 42 - VeChain
Set-TextValue "D42" "0.0227"
Set-TextValue "E42" "  +3.91%  "

# Row 43 - Maker
Set-TextValue "D43" "1.527.94"
Set-TextValue "E43" "  -0.60%  "

# Row 44 - TrustWalletToken
Set-TextValue "E44" "  +8.40%  "

# Row 45 - HuobiToken
Set-TextValue "E45" "  -0.64%  "

# Row 46 - FraxShare
Set-TextValue "D46" "7.80"
Set-TextValue "E46" "  +0.02%  "

# Row 47 - Cronos
Set-TextValue "D47" "0.0917"
Set-TextValue "E47" "  -2.42%  "

# Row 48 - ARBITRUM
Set-TextValue "E48" "  +4.67%  "

# Row 49 - FTXToken
Set-TextValue "D49" "4.15"
Set-TextValue "E49" "  +0.01%  "

# Row 50 - MXToken
Set-TextValue "E50" "  -0.49%  "

# Row 51 - RocketPoolETH
Set-TextValue "D51" "2.303.52"
Set-TextValue "E51" "  +0.82%  "
